$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Sheet 2")

# Swap the numeric values that live in D9 and E9
$ws.Range("D9").Value = 235
$ws.Range("E9").Value = 234

# Make the sheet active and move the selection to E9
$ws.Activate()
$ws.Range("E9").Select()
